# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# 1) Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 9 de Mayo de 2020 a las 07:04"

# 2) Row 4 - Estados Unidos: refreshed totals
$ws.Range("B4").Value = 1322163
$ws.Range("C4").Value = 378
$ws.Range("D4").Value = 223749
$ws.Range("E4").Value = 1019798
$ws.Range("F4").Value = 16978
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 78616

# 3) Row 69 - Tailandia: refreshed totals
$ws.Range("B69").Value = 3004
$ws.Range("C69").Value = 4
$ws.Range("D69").Value = 2787
$ws.Range("E69").Value = 161
$ws.Range("F69").Value = 61
$ws.Range("G69").Value = 1
$ws.Range("H69").Value = 56

# 4) A new country "El Salvador" is inserted into the ranking right after
#    "Niger" (row 106), which shifts Costa Rica / Principado de Andorra /
#    Maldivas / Burkina Faso down one row each (rows 107-110), and the
#    data that used to belong to Burkina Faso's row (110) is replaced by
#    the data that used to sit at row 109 (what was previously beneath it).

# Row 106: now El Salvador, with freshly reported stats
$ws.Range("A106").Value = "El Salvador"
$ws.Range("B106").Value = 784
$ws.Range("C106").Value = 42
$ws.Range("D106").Value = 276
$ws.Range("E106").Value = 492
$ws.Range("F106").Value = 4
$ws.Range("G106").Value = 0
$ws.Range("H106").Value = 16

# Row 107: now Costa Rica (old row 106 values)
$ws.Range("A107").Value = "Costa Rica"
$ws.Range("B107").Value = 773
$ws.Range("C107").Value = 0
$ws.Range("D107").Value = 461
$ws.Range("E107").Value = 306
$ws.Range("F107").Value = 6
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = 6

# Row 108: now Principado de Andorra (old row 107 values)
$ws.Range("A108").Value = "Principado de Andorra"
$ws.Range("B108").Value = 752
$ws.Range("C108").Value = 0
$ws.Range("D108").Value = 537
$ws.Range("E108").Value = 168
$ws.Range("F108").Value = 14
$ws.Range("G108").Value = 0
$ws.Range("H108").Value = 47

# Row 109: now Maldivas (old row 108 values)
$ws.Range("A109").Value = "Maldivas"
$ws.Range("B109").Value = 744
$ws.Range("C109").Value = 0
$ws.Range("D109").Value = 20
$ws.Range("E109").Value = 721
$ws.Range("F109").Value = 2
$ws.Range("G109").Value = 0
$ws.Range("H109").Value = 3

# Row 110: now Burkina Faso (old row 109 values)
$ws.Range("A110").Value = "Burkina Faso"
$ws.Range("B110").Value = 744
$ws.Range("C110").Value = 0
$ws.Range("D110").Value = 566
$ws.Range("E110").Value = 130
$ws.Range("F110").Value = 0
$ws.Range("G110").Value = 0
$ws.Range("H110").Value = 48

# 5) Belice and Nueva Caledonia swap places (rows 192 / 193)
$ws.Range("A192").Value = "Nueva Caledonia"
$ws.Range("D192").Value = 18
$ws.Range("H192").Value = 0

$ws.Range("A193").Value = "Belice"
$ws.Range("D193").Value = 16
$ws.Range("H193").Value = 2

# 6) Row 212 - Islas Virgenes Britanicas: refreshed totals
$ws.Range("D212").Value = 4
$ws.Range("E212").Value = 2
